# Commit: "Fixed #295 Add the version of M2Doc in the template custom
# properties."
#
# For THIS fixture (invalidConditionInElseif-template.docx) the supplied
# canonical-OOXML diff contains no real content change at all. Every one
# of its ~45 changed lines is the very same element, with the very same
# attribute names and values, just re-serialized with the attributes in
# (mostly) alphabetical order instead of their original order, e.g.:
#
#   <w:tab w:val="left" w:pos="3119"/>  ->  <w:tab w:pos="3119" w:val="left"/>
#   <w:pgSz w:w="11906" w:h="16838"/>   ->  <w:pgSz w:h="16838" w:w="11906"/>
#   <w:style w:type="paragraph" w:default="1" w:styleId="Normal">
#                                       ->  <w:style w:default="1" w:styleId="Normal" w:type="paragraph">
#
# and likewise for the document's xmlns declarations, the docDefaults
# rFonts/lang, every w:latentStyles/w:lsdException entry, and every
# w:style tag. No element/attribute/value/text is added, removed or
# changed anywhere - the attribute sets are identical, only the byte
# order of attributes within a start tag differs. That is a leftover of
# whatever tool re-saved the fixture for the upstream commit (the actual
# #295 feature - stamping the M2Doc version into the template's custom
# document properties - lives in other files from that commit, not in
# this template). Word's object model has no "attribute order" knob to
# begin with, so there is nothing for a COM script to do to reproduce
# that byte shuffle - the paragraphs/runs/tab stops/section page setup/
# fonts/language/styles this diff touches are already exactly right.
#
# Verify (read-only - no setters, so nothing gets marked dirty/rewritten
# and the package stays byte-for-byte what it already correctly is).

$d = $word.ActiveDocument

$tabHits = 0
foreach ($p in $d.Paragraphs) {
    $tabs = $p.Range.ParagraphFormat.TabStops
    for ($i = 1; $i -le $tabs.Count; $i++) {
        if ([math]::Round($tabs.Item($i).Position * 20) -eq 3119) {
            $tabHits = $tabHits + 1
        }
    }
}
Write-Output "left tab stops at 3119 twips: $tabHits"

$ps = $d.Sections.Item(1).PageSetup
$w = [math]::Round($ps.PageWidth * 20)
$h = [math]::Round($ps.PageHeight * 20)
Write-Output "page size (twips): $w x $h"

$normal = $d.Styles("Normal")
Write-Output "Normal style font: $($normal.Font.NameAscii) $($normal.Font.Size)pt"
